$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Spp1"
$ws.Range("C2").Value = "Itgb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 209.0063303333334
$ws.Range("H2").Value = 627.018991
$ws.Range("I2").Value = 0.6751081226665357
$ws.Range("J2").Value = 0.6751081226665357
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 112.513392
$ws.Range("N2").Value = 337.540176
$ws.Range("O2").Value = 0.3275312977368564
$ws.Range("P2").Value = 0.3275312977368564
$ws.Range("Q2").Value = 23516.01117527582
$ws.Range("R2").Value = 211644.1005774824
$ws.Range("S2").Value = 0.2211190395296633
$ws.Range("T2").Value = 0.2211190395296633

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Spp1"
$ws.Range("C3").Value = "Itgb1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 209.0063303333334
$ws.Range("H3").Value = 627.018991
$ws.Range("I3").Value = 0.6751081226665357
$ws.Range("J3").Value = 0.6751081226665357
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 106.314466
$ws.Range("N3").Value = 318.943398
$ws.Range("O3").Value = 0.3094859589441663
$ws.Range("P3").Value = 0.3094859589441664
$ws.Range("Q3").Value = 22220.39640000794
$ws.Range("R3").Value = 199983.5676000714
$ws.Range("S3").Value = 0.2089364847344487
$ws.Range("T3").Value = 0.2089364847344487

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Spp1"
$ws.Range("C4").Value = "Itgb1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 209.0063303333334
$ws.Range("H4").Value = 627.018991
$ws.Range("I4").Value = 0.6751081226665357
$ws.Range("J4").Value = 0.6751081226665357
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 124.6916553333333
$ws.Range("N4").Value = 374.074966
$ws.Range("O4").Value = 0.3629827433189773
$ws.Range("P4").Value = 0.3629827433189773
$ws.Range("Q4").Value = 26061.34530440882
$ws.Range("R4").Value = 234552.1077396793
$ws.Range("S4").Value = 0.2450525984024238
$ws.Range("T4").Value = 0.2450525984024238

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Spp1"
$ws.Range("C5").Value = "Itgb1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.401741666666666
$ws.Range("H5").Value = 4.205225
$ws.Range("I5").Value = 0.004527744128790482
$ws.Range("J5").Value = 0.004527744128790482
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 112.513392
$ws.Range("N5").Value = 337.540176
$ws.Range("O5").Value = 0.3275312977368564
$ws.Range("P5").Value = 0.3275312977368564
$ws.Range("Q5").Value = 157.7147096244
$ws.Range("R5").Value = 1419.4323866196
$ws.Range("S5").Value = 0.001482977910323179
$ws.Range("T5").Value = 0.001482977910323179

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Spp1"
$ws.Range("C6").Value = "Itgb1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.401741666666666
$ws.Range("H6").Value = 4.205225
$ws.Range("I6").Value = 0.004527744128790482
$ws.Range("J6").Value = 0.004527744128790482
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 106.314466
$ws.Range("N6").Value = 318.943398
$ws.Range("O6").Value = 0.3094859589441663
$ws.Range("P6").Value = 0.3094859589441664
$ws.Range("Q6").Value = 149.0254167616166
$ws.Range("R6").Value = 1341.22875085455
$ws.Range("S6").Value = 0.001401273233552541
$ws.Range("T6").Value = 0.001401273233552541

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Spp1"
$ws.Range("C7").Value = "Itgb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 1.401741666666666
$ws.Range("H7").Value = 4.205225
$ws.Range("I7").Value = 0.004527744128790482
$ws.Range("J7").Value = 0.004527744128790482
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 124.6916553333333
$ws.Range("N7").Value = 374.074966
$ws.Range("O7").Value = 0.3629827433189773
$ws.Range("P7").Value = 0.3629827433189773
$ws.Range("Q7").Value = 174.7854887663722
$ws.Range("R7").Value = 1573.06939889735
$ws.Range("S7").Value = 0.001643492984914762
$ws.Range("T7").Value = 0.001643492984914762

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Spp1"
$ws.Range("C8").Value = "Itgb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 99.18134533333334
$ws.Range("H8").Value = 297.544036
$ws.Range("I8").Value = 0.3203641332046738
$ws.Range("J8").Value = 0.3203641332046737
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 112.513392
$ws.Range("N8").Value = 337.540176
$ws.Range("O8").Value = 0.3275312977368564
$ws.Range("P8").Value = 0.3275312977368564
$ws.Range("Q8").Value = 11159.2295865767
$ws.Range("R8").Value = 100433.0662791903
$ws.Range("S8").Value = 0.1049292802968699
$ws.Range("T8").Value = 0.1049292802968699

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Spp1"
$ws.Range("C9").Value = "Itgb1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 99.18134533333334
$ws.Range("H9").Value = 297.544036
$ws.Range("I9").Value = 0.3203641332046738
$ws.Range("J9").Value = 0.3203641332046737
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 106.314466
$ws.Range("N9").Value = 318.943398
$ws.Range("O9").Value = 0.3094859589441663
$ws.Range("P9").Value = 0.3094859589441664
$ws.Range("Q9").Value = 10544.41176627492
$ws.Range("R9").Value = 94899.70589647433
$ws.Range("S9").Value = 0.0991482009761651
$ws.Range("T9").Value = 0.0991482009761651

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Spp1"
$ws.Range("C10").Value = "Itgb1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 99.18134533333334
$ws.Range("H10").Value = 297.544036
$ws.Range("I10").Value = 0.3203641332046738
$ws.Range("J10").Value = 0.3203641332046737
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 124.6916553333333
$ws.Range("N10").Value = 374.074966
$ws.Range("O10").Value = 0.3629827433189773
$ws.Range("P10").Value = 0.3629827433189773
$ws.Range("Q10").Value = 12367.08612780031
$ws.Range("R10").Value = 111303.7751502028
$ws.Range("S10").Value = 0.1162866519316387
$ws.Range("T10").Value = 0.1162866519316387

